# anulação da regra que removia anos com dados faltantes e adição de
# coluna que indica se há dados faltantes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F ("Faltam dados para todos os Estados"), styled
# like the other header cells (bold, thin border, centered/top aligned).
$ws.Range("F1").Value = "Faltam dados para todos os Estados"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160

# Full replacement data set: rows now span 2015-2024 for Brasil, Nordeste
# and Sergipe (previously only 2019-2024), plus the new "Faltam dados"
# boolean flag column (F) that marks years where data was missing for all
# states (the rule that used to drop those years was removed).
$data = @(
    @(2, "Brasil", "01/01/2015", "Feminicídio", 0.618096459842882, $null, 1),
    @(3, "Brasil", "01/01/2016", "Feminicídio", 0.863032297756537, $null, 1),
    @(4, "Brasil", "01/01/2017", "Feminicídio", 1.082538168435736, $null, 1),
    @(5, "Brasil", "01/01/2018", "Feminicídio", 1.228029758579693, $null, 1),
    @(6, "Brasil", "01/01/2019", "Feminicídio", 1.491603669709312, $null, 0),
    @(7, "Brasil", "01/01/2020", "Feminicídio", 1.53032839946819, $null, 0),
    @(8, "Brasil", "01/01/2021", "Feminicídio", 1.569654974814453, $null, 0),
    @(9, "Brasil", "01/01/2022", "Feminicídio", 1.532101471544391, $null, 0),
    @(10, "Brasil", "01/01/2023", "Feminicídio", 1.503723149276654, $null, 0),
    @(11, "Brasil", "01/01/2024", "Feminicídio", 1.443586697474013, $null, 0),
    @(12, "Nordeste", "01/01/2015", "Feminicídio", 0.6635307538337432, $null, 1),
    @(13, "Nordeste", "01/01/2016", "Feminicídio", 0.9129393838946243, $null, 1),
    @(14, "Nordeste", "01/01/2017", "Feminicídio", 1.239969626186528, $null, 1),
    @(15, "Nordeste", "01/01/2018", "Feminicídio", 1.332959341705364, $null, 1),
    @(16, "Nordeste", "01/01/2019", "Feminicídio", 1.497286779739304, $null, 0),
    @(17, "Nordeste", "01/01/2020", "Feminicídio", 1.421470954921448, $null, 0),
    @(18, "Nordeste", "01/01/2021", "Feminicídio", 1.450454625600147, $null, 0),
    @(19, "Nordeste", "01/01/2022", "Feminicídio", 1.326824150475039, $null, 0),
    @(20, "Nordeste", "01/01/2023", "Feminicídio", 1.33400727814508, $null, 0),
    @(21, "Nordeste", "01/01/2024", "Feminicídio", 1.31935864980953, $null, 0),
    @(22, "Sergipe", "01/01/2015", "Feminicídio", 0, 20.5, 1),
    @(23, "Sergipe", "01/01/2016", "Feminicídio", 0, 23, 1),
    @(24, "Sergipe", "01/01/2017", "Feminicídio", 1.715876576997817, 5, 1),
    @(25, "Sergipe", "01/01/2018", "Feminicídio", 1.359138849624878, 13, 1),
    @(26, "Sergipe", "01/01/2019", "Feminicídio", 1.766753237575308, 7, 0),
    @(27, "Sergipe", "01/01/2020", "Feminicídio", 1.166870869068754, 20, 0),
    @(28, "Sergipe", "01/01/2021", "Feminicídio", 1.651913700724447, 9, 0),
    @(29, "Sergipe", "01/01/2022", "Feminicídio", 1.555611048113412, 12, 0),
    @(30, "Sergipe", "01/01/2023", "Feminicídio", 1.298940632730219, 19, 0),
    @(31, "Sergipe", "01/01/2024", "Feminicídio", 0.8052392083854389, 25, 0)
)

# Column B holds dates written as literal text (e.g. "01/01/2015"), not
# real date values. Force text format first so Excel doesn't silently
# convert the strings into date serials, then strip the format change
# back off (copying the neighboring default-styled cell's style) so the
# cells end up with no special number formatting, same as the source.
$ws.Range("B2:B31").NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    if ($row[5] -eq $null) {
        $ws.Cells.Item($r, 5).Value = ""
    } else {
        $ws.Cells.Item($r, 5).Value = $row[5]
    }
    $ws.Cells.Item($r, 6).Value = [bool]$row[6]
}

$ws.Range("B2:B31").Style = $ws.Range("A2").Style

Write-Host "Applied g19.7 update: expanded years + Faltam dados column"
